$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: set a plain text value into a cell without Excel's automatic
# date/number inference turning date-like strings (e.g. "2012.6.10")
# into date serials. We build the literal as a formula string and then
# replace the formula with its computed value (paste values), which
# yields a plain shared-string cell with no special number format.
function Set-PlainText($addr, $text) {
    $cell = $ws.Range($addr)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null # xlPasteValues
}

# Populate the new shared strings in the same order the original author
# entered them, so that the shared-string table indices line up exactly.
Set-PlainText "A30" "2012.6.17"
Set-PlainText "B30" "根据设计，代码修改实现20关"
Set-PlainText "B28" "20关卡设计"
Set-PlainText "B29" "小组会议"
Set-PlainText "B27" "1.小鸟控制修改：改为手势控制；2消分判断修改：仓库最后3个一样才能触发消分。"
Set-PlainText "C27" "消分时现在是从左到右遍历，需要改为从右向左遍历"
Set-PlainText "A29" "2012.6.15"
Set-PlainText "A28" "2012.6.14"
Set-PlainText "C26" "碰撞后糖果运行轨迹还需要更好的算法"
Set-PlainText "C30" "后续可以加入包装的糖果"
Set-PlainText "A27" "2012.6.10"
Set-PlainText "A31" "2012.6.21"
Set-PlainText "B31" "加入包装糖果，修改消分时从右往左遍历"

# --- Formatting / remaining (non-string) values ---
$ws.Range("B27").WrapText = $true
$ws.Range("B28").WrapText = $true
$ws.Range("B29").WrapText = $true
$ws.Range("B31").WrapText = $true
$ws.Rows.Item(27).RowHeight = 27
$ws.Range("D27").Value = 3
$ws.Range("D28").Value = 2
$ws.Range("D30").Value = 5
$ws.Range("D31").Value = 3.5

# --- Update the selection to match the edited workbook ---
$ws.Range("C36").Select()
